$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "NegativeTests" worksheet between "DataSet" and
#    "DataSetInteractionPages".
# ---------------------------------------------------------------------------
$dataSetSheet = $wb.Worksheets.Item("DataSet")
$negSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSetSheet)
$negSheet.Name = "NegativeTests"

# Give the whole used range a text number format (matches the rest of the
# workbook, which stores everything as text - numFmtId 49).
$negSheet.Range("A1:J5").NumberFormat = "@"

# ---- Header row ------------------------------------------------------------
$negSheet.Range("A1").Value = "Key"
$negSheet.Range("B1").Value = "Email"
$negSheet.Range("C1").Value = "Fullname"
$negSheet.Range("D1").Value = "Password"
$negSheet.Range("E1").Value = "ConfirmPassword"
$negSheet.Range("F1").Value = "Effect"
$negSheet.Range("G1").Value = "Effect2"
$negSheet.Range("H1").Value = "Asserter"
$negSheet.Range("I1").Value = "Asserter2"
$negSheet.Range("J1").Value = "Status"

# ---- Row 2: Register_Without_Email -----------------------------------------
$negSheet.Range("A2").Value = "Register_Without_Email"
$negSheet.Range("C2").Value = "test"
$negSheet.Range("D2").Value = "test"
$negSheet.Range("E2").Value = "test"
$negSheet.Range("F2").Value = "The Email field is required."
$negSheet.Range("H2").Value = "AssertEmailErrorMessageExists"
$negSheet.Range("J2").Value = "'passed"

# ---- Row 3: Register_Without_FullName --------------------------------------
$negSheet.Range("A3").Value = "Register_Without_FullName"
$negSheet.Range("B3").Value = "test@test.com"
$negSheet.Range("D3").Value = "test"
$negSheet.Range("E3").Value = "test"
$negSheet.Range("F3").Value = "The Full Name field is required."
$negSheet.Range("H3").Value = "AssertFullNameErrorMessageExists"
$negSheet.Range("J3").Value = "'passed"
$negSheet.Hyperlinks.Add($negSheet.Range("B3"), "mailto:test@test.com") | Out-Null

# ---- Row 4: Register_Without_Password ---------------------------------------
$negSheet.Range("A4").Value = "Register_Without_Password"
$negSheet.Range("B4").Value = "test@test.com"
$negSheet.Range("C4").Value = "test"
$negSheet.Range("D4").Value = " "
$negSheet.Range("E4").Value = "test"
$negSheet.Range("F4").Value = "The Password field is required."
$negSheet.Range("G4").Value = "The password and confirmation password do not match."
$negSheet.Range("H4").Value = "AssertPasswordErrorMessageExists"
$negSheet.Range("I4").Value = "AssertPasswordMissmatchErrorMessageExists"
$negSheet.Range("J4").Value = "'passed"
$negSheet.Hyperlinks.Add($negSheet.Range("B4"), "mailto:test@test.com") | Out-Null

# ---- Row 5: Register_Without_ConfirmPassword --------------------------------
$negSheet.Range("A5").Value = "Register_Without_ConfirmPassword"
$negSheet.Range("B5").Value = "test@test.com"
$negSheet.Range("C5").Value = "test"
$negSheet.Range("D5").Value = "test"
$negSheet.Range("F5").Value = "The password and confirmation password do not match."
$negSheet.Range("H5").Value = "AssertPasswordMissmatchErrorMessageExists2"
$negSheet.Hyperlinks.Add($negSheet.Range("B5"), "mailto:test@test.com") | Out-Null

# Page setup matching the rest of the workbook's worksheets.
$negSheet.PageSetup.PaperSize = 9
$negSheet.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 2. DataSet sheet: one of the existing test results actually flips from
#    "failed" to "passed" (row 4 / column P), and the status cells swap from
#    inline strings to shared strings as part of the resave.  Prefixing the
#    value with an apostrophe preserves each cell's existing "quote
#    prefix" number format (so the style index is left untouched).
# ---------------------------------------------------------------------------
$dataSet = $wb.Worksheets.Item("DataSet")
$dataSet.Range("P4").Value = "'passed"
$dataSet.Range("P7").Value = "'passed"
$dataSet.Range("P10").Value = "'failed"
$dataSet.Range("P15").Value = "'failed"
$dataSet.Range("P18").Value = "'passed"
$dataSet.Range("P22").Value = "passed"

# ---------------------------------------------------------------------------
# 3. Make the new sheet the active / selected tab, matching the author's
#    final view state.
# ---------------------------------------------------------------------------
$negSheet.Range("H5").Select()
$negSheet.Activate()
